# "Processed tabs and figs"
# - Row 7 (A7) / Row 8 (A8): swap the sample labels ("纯水" <-> the
#   "0.05%CuSO...(aq)" label) that had been entered in the wrong rows.
# - Row 6 (B6): the placeholder "--" reading is removed (cell content
#   cleared, keeping its existing border/quote-prefix style).
# - Scroll the sheet view over a column and move the active selection
#   (cosmetic view-state housekeeping that goes with the data fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6 no longer holds the "--" stand-in reading.
$ws.Range("B6").ClearContents()

# A7 / A8 had their sample-name strings swapped.
$ws.Range("A7").Value = '纯水'
$ws.Range("A8").Value = '0.05\%CuSO\textsubscript{9}(aq)'

# Scroll the view right one column and move the selection, matching the
# saved view state (topLeftCell=B1, selection=N24).
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$ws.Range("N24").Select()
